$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 24.916566
$ws.Cells.Item(2, 8).Value = 74.749698
$ws.Cells.Item(2, 9).Value = 0.459912889255076
$ws.Cells.Item(2, 10).Value = 0.459912889255076
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 33.380049
$ws.Cells.Item(2, 14).Value = 100.140147
$ws.Cells.Item(2, 15).Value = 0.3891462059670435
$ws.Cells.Item(2, 16).Value = 0.3891462059670435
$ws.Cells.Item(2, 17).Value = 831.716193991734
$ws.Cells.Item(2, 18).Value = 7485.445745925605
$ws.Cells.Item(2, 19).Value = 0.1789733559289539
$ws.Cells.Item(2, 20).Value = 0.1789733559289539

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 24.916566
$ws.Cells.Item(3, 8).Value = 74.749698
$ws.Cells.Item(3, 9).Value = 0.459912889255076
$ws.Cells.Item(3, 10).Value = 0.459912889255076
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 33.85786133333334
$ws.Cells.Item(3, 14).Value = 101.573584
$ws.Cells.Item(3, 15).Value = 0.3947165649764305
$ws.Cells.Item(3, 16).Value = 0.3947165649764305
$ws.Cells.Item(3, 17).Value = 843.6216365308481
$ws.Cells.Item(3, 18).Value = 7592.594728777633
$ws.Cells.Item(3, 19).Value = 0.1815352358351491
$ws.Cells.Item(3, 20).Value = 0.1815352358351491

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 24.916566
$ws.Cells.Item(4, 8).Value = 74.749698
$ws.Cells.Item(4, 9).Value = 0.459912889255076
$ws.Cells.Item(4, 10).Value = 0.459912889255076
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 18.53974466666667
$ws.Cells.Item(4, 14).Value = 55.61923400000001
$ws.Cells.Item(4, 15).Value = 0.2161372290565261
$ws.Cells.Item(4, 16).Value = 0.2161372290565261
$ws.Cells.Item(4, 17).Value = 461.946771610148
$ws.Cells.Item(4, 18).Value = 4157.520944491333
$ws.Cells.Item(4, 19).Value = 0.09940429749097307
$ws.Cells.Item(4, 20).Value = 0.09940429749097307

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 19.60300333333333
$ws.Cells.Item(5, 8).Value = 58.80901
$ws.Cells.Item(5, 9).Value = 0.3618345281251927
$ws.Cells.Item(5, 10).Value = 0.3618345281251927
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 33.380049
$ws.Cells.Item(5, 14).Value = 100.140147
$ws.Cells.Item(5, 15).Value = 0.3891462059670435
$ws.Cells.Item(5, 16).Value = 0.3891462059670435
$ws.Cells.Item(5, 17).Value = 654.34921181383
$ws.Cells.Item(5, 18).Value = 5889.14290632447
$ws.Cells.Item(5, 19).Value = 0.1408065338077942
$ws.Cells.Item(5, 20).Value = 0.1408065338077942

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 19.60300333333333
$ws.Cells.Item(6, 8).Value = 58.80901
$ws.Cells.Item(6, 9).Value = 0.3618345281251927
$ws.Cells.Item(6, 10).Value = 0.3618345281251927
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 33.85786133333334
$ws.Cells.Item(6, 14).Value = 101.573584
$ws.Cells.Item(6, 15).Value = 0.3947165649764305
$ws.Cells.Item(6, 16).Value = 0.3947165649764305
$ws.Cells.Item(6, 17).Value = 663.7157685768713
$ws.Cells.Item(6, 18).Value = 5973.441917191841
$ws.Cells.Item(6, 19).Value = 0.1428220820314437
$ws.Cells.Item(6, 20).Value = 0.1428220820314437

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 19.60300333333333
$ws.Cells.Item(7, 8).Value = 58.80901
$ws.Cells.Item(7, 9).Value = 0.3618345281251927
$ws.Cells.Item(7, 10).Value = 0.3618345281251927
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 18.53974466666667
$ws.Cells.Item(7, 14).Value = 55.61923400000001
$ws.Cells.Item(7, 15).Value = 0.2161372290565261
$ws.Cells.Item(7, 16).Value = 0.2161372290565261
$ws.Cells.Item(7, 17).Value = 363.4346764998156
$ws.Cells.Item(7, 18).Value = 3270.91208849834
$ws.Cells.Item(7, 19).Value = 0.07820591228595479
$ws.Cells.Item(7, 20).Value = 0.07820591228595479

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 9.657138000000002
$ws.Cells.Item(8, 8).Value = 28.971414
$ws.Cells.Item(8, 9).Value = 0.1782525826197313
$ws.Cells.Item(8, 10).Value = 0.1782525826197313
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 33.380049
$ws.Cells.Item(8, 14).Value = 100.140147
$ws.Cells.Item(8, 15).Value = 0.3891462059670435
$ws.Cells.Item(8, 16).Value = 0.3891462059670435
$ws.Cells.Item(8, 17).Value = 322.355739639762
$ws.Cells.Item(8, 18).Value = 2901.201656757858
$ws.Cells.Item(8, 19).Value = 0.06936631623029538
$ws.Cells.Item(8, 20).Value = 0.06936631623029538

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 9.657138000000002
$ws.Cells.Item(9, 8).Value = 28.971414
$ws.Cells.Item(9, 9).Value = 0.1782525826197313
$ws.Cells.Item(9, 10).Value = 0.1782525826197313
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 33.85786133333334
$ws.Cells.Item(9, 14).Value = 101.573584
$ws.Cells.Item(9, 15).Value = 0.3947165649764305
$ws.Cells.Item(9, 16).Value = 0.3947165649764305
$ws.Cells.Item(9, 17).Value = 326.9700392808641
$ws.Cells.Item(9, 18).Value = 2942.730353527777
$ws.Cells.Item(9, 19).Value = 0.07035924710983771
$ws.Cells.Item(9, 20).Value = 0.0703592471098377

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 9.657138000000002
$ws.Cells.Item(10, 8).Value = 28.971414
$ws.Cells.Item(10, 9).Value = 0.1782525826197313
$ws.Cells.Item(10, 10).Value = 0.1782525826197313
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 18.53974466666667
$ws.Cells.Item(10, 14).Value = 55.61923400000001
$ws.Cells.Item(10, 15).Value = 0.2161372290565261
$ws.Cells.Item(10, 16).Value = 0.2161372290565261
$ws.Cells.Item(10, 17).Value = 179.040872730764
$ws.Cells.Item(10, 18).Value = 1611.367854576876
$ws.Cells.Item(10, 19).Value = 0.03852701927959819
$ws.Cells.Item(10, 20).Value = 0.03852701927959819

Write-Host "Updated rows 2-10 columns E:T"